# Apply updated crypto price/volume figures to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D occasionally hold values that look like plain numbers
# (e.g. "106.00", "0.542"). The sheet stores these as text so that formatting
# such as trailing zeros is preserved, so we prefix them with a leading
# apostrophe (Excel's standard "treat as text" marker) before assigning them.

$ws.Range('D2').Value = '47.828.45'
$ws.Range('E2').Value = '  +5.90%  '
$ws.Range('D3').Value = '2.510.14'
$ws.Range('E3').Value = '  +3.39%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''324.25'
$ws.Range('E5').Value = '  +2.42%  '
$ws.Range('D6').Value = '''106.00'
$ws.Range('E6').Value = '  +3.32%  '
$ws.Range('E7').Value = '  +1.65%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '''0.542'
$ws.Range('E9').Value = '  +3.13%  '
$ws.Range('D10').Value = '''38.02'
$ws.Range('E10').Value = '  +7.17%  '
$ws.Range('E11').Value = '  +1.74%  '
$ws.Range('D12').Value = '''0.124'
$ws.Range('E12').Value = '  +1.03%  '
$ws.Range('D13').Value = '''18.42'
$ws.Range('E13').Value = '  +1.55%  '
$ws.Range('E14').Value = '  +2.28%  '
$ws.Range('D15').Value = '2.904.01'
$ws.Range('E15').Value = '  +3.38%  '
$ws.Range('D16').Value = '2.521.53'
$ws.Range('E16').Value = '  +3.41%  '
$ws.Range('E17').Value = '  +1.12%  '
$ws.Range('D18').Value = '47.693.22'
$ws.Range('E18').Value = '  +5.79%  '
$ws.Range('D19').Value = '''12.76'
$ws.Range('E19').Value = '  +3.99%  '
$ws.Range('E20').Value = '  +3.48%  '
$ws.Range('E21').Value = '  +2.00%  '
$ws.Range('D22').Value = '''70.84'
$ws.Range('D23').Value = '''251.54'
$ws.Range('E23').Value = '  +3.17%  '
$ws.Range('E24').Value = '  +6.55%  '
$ws.Range('E25').Value = '  +3.10%  '
$ws.Range('D26').Value = '''26.35'
$ws.Range('E26').Value = '  +3.32%  '
$ws.Range('E27').Value = '  -0.11%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = '''2.28'
$ws.Range('E28').Value = '  +10.49%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').Value = '''10.04'
$ws.Range('E29').Value = '  +4.93%  '
$ws.Range('D30').Value = '''35.28'
$ws.Range('E31').Value = '  +9.62%  '
$ws.Range('E32').Value = '  +0.49%  '
$ws.Range('D33').Value = '''20.05'
$ws.Range('E33').Value = '  -1.35%  '
$ws.Range('E34').Value = '  +3.22%  '
$ws.Range('D35').Value = '''0.0784'
$ws.Range('E35').Value = '  +2.74%  '
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('E38').Value = '  +4.73%  '
$ws.Range('E39').Value = '  +4.91%  '
$ws.Range('E40').Value = '  +2.13%  '
$ws.Range('E41').Value = '  +2.23%  '
$ws.Range('D42').Value = '''121.61'
$ws.Range('E42').Value = '  -2.24%  '
$ws.Range('D43').Value = '''21.22'
$ws.Range('E43').Value = '  +3.16%  '
$ws.Range('D44').Value = '''0.0299'
$ws.Range('E44').Value = '  +3.73%  '
$ws.Range('D45').Value = '1.971.83'
$ws.Range('D46').Value = '''3.01'
$ws.Range('E46').Value = '  +3.18%  '
$ws.Range('D47').Value = '''2.10'
$ws.Range('E47').Value = '  -0.34%  '
$ws.Range('E48').Value = '  +0.67%  '
$ws.Range('D49').Value = '''9.24'
$ws.Range('E49').Value = '  +0.49%  '
$ws.Range('E50').Value = '  +14.05%  '
$ws.Range('D51').Value = '''79.25'
$ws.Range('E51').Value = '  +3.59%  '
